$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Columns D/E hold plain text in the source data (inline strings),
# e.g. "1.00", "0.0407", "582.30". Assigning such strings straight to
# Range.Value would let Excel auto-coerce them into numbers (dropping
# trailing zeros / switching to scientific notation), so for any new
# value that looks numeric we force the cell to Text format first.

# Row 2
$ws.Range("D2").Value = "67.823.45"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").Value = "3.318.59"
$ws.Range("E3").Value = "  -1.70%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.30"
$ws.Range("E5").Value = "  -1.93%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.61"
$ws.Range("E6").Value = "  -6.62%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("E8").Value = "  -2.34%  "

# Row 9
$ws.Range("D9").Value = "3.316.65"
$ws.Range("E9").Value = "  -1.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  -4.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.576"
$ws.Range("E11").Value = "  -2.47%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.36"
$ws.Range("E12").Value = "  -4.81%  "

# Row 13
$ws.Range("E13").Value = "  -2.81%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "665.53"
$ws.Range("E14").Value = "  +3.70%  "

# Row 15
$ws.Range("D15").Value = "3.857.56"
$ws.Range("E15").Value = "  -1.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.38"
$ws.Range("E16").Value = "  -2.98%  "

# Row 17
$ws.Range("D17").Value = "67.887.98"

# Row 18
$ws.Range("E18").Value = "  -0.97%  "

# Row 19
$ws.Range("D19").Value = "3.320.56"
$ws.Range("E19").Value = "  -1.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.45"
$ws.Range("E20").Value = "  -3.59%  "

# Row 21
$ws.Range("E21").Value = "  -2.58%  "

# Row 22
$ws.Range("E22").Value = "  -2.85%  "

# Row 23
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.35"
$ws.Range("E23").Value = "  +4.60%  "

# Row 24
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.99"
$ws.Range("E24").Value = "  -5.96%  "

# Row 25
$ws.Range("E25").Value = "  -2.59%  "

# Row 26
$ws.Range("E26").Value = "  -5.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -6.57%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.26"
$ws.Range("E28").Value = "  -5.80%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.72"
$ws.Range("E29").Value = "  +3.31%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.40"
$ws.Range("E30").Value = "  -3.73%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.31"
$ws.Range("E31").Value = "  +4.56%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "585.38"
$ws.Range("E32").Value = "  -4.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.96"
$ws.Range("E33").Value = "  -1.62%  "

# Row 34
$ws.Range("E34").Value = "  -2.66%  "

# Row 35
$ws.Range("D35").Value = "3.724.97"
$ws.Range("E35").Value = "  -7.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.28"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  -13.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.131"
$ws.Range("E39").Value = "  -0.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.51"
$ws.Range("E40").Value = "  -3.90%  "

# Row 41
$ws.Range("E41").Value = "  -6.97%  "

# Row 42
$ws.Range("E42").Value = "  -5.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.332"
$ws.Range("E43").Value = "  -3.57%  "

# Row 44
$ws.Range("D44").Value = "0.0₃0663"
$ws.Range("E44").Value = "  -6.29%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -5.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  -4.14%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.60"
$ws.Range("E47").Value = "  -0.29%  "

# Row 48
$ws.Range("E48").Value = "  -2.19%  "

# Row 49
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("E50").Value = "  -2.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.52"
$ws.Range("E51").Value = "  -0.88%  "
